$d = $word.ActiveDocument
$found1 = $d.Content.Find.Execute("{Eventueel kunnen stagebedrijfen hier ook informatie toevoegen}", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
Write-Output $found1
$found2 = $d.Content.Find.Execute(". Le", $true, $false, $false, $false, $false, $true, 1, $false, "Le", 2)
Write-Output $found2
